# "added error bars for x data"
#
# The Part1 worksheet's two side-by-side measurement tables each get a new
# "error (cm)" column inserted right after their "deltah(cm)" column
# (i.e. before the existing "V (mV)" column), filled with a constant
# error value of 1. The line chart (which lives on Sheet1 but plots data
# from Part1) has its voltage series updated to follow the column that
# now holds the "V (mV)" values (shifted from column B to column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Part1")

# --- Insert the two new "error (cm)" columns -----------------------------
# First table: A=deltah(cm), B=V(mV), C=error(mV)  -> insert blank col at B
$ws.Range("B1").EntireColumn.Insert()
# Second table now starts at E (was D) -> insert blank col before its V(mV)
# column, i.e. before what is now F (was E)
$ws.Range("F1").EntireColumn.Insert()

# --- Fill in the new column headers --------------------------------------
$ws.Range("B1").Value = "error (cm)"
$ws.Range("F1").Value = "error (cm)"

# --- Fill in the new column values (constant error of 1) -----------------
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("B$row").Value = 1
    $ws.Range("F$row").Value = 1
}

# --- Update the selection shown in the sheet view -------------------------
[void]$ws.Range("B1").Select()

# --- Fix up the chart: the voltage series now lives in column C ----------
$wsChart = $wb.Worksheets.Item("Sheet1")
$chartObj = $wsChart.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = '=SERIES(,Part1!$A$1:$A$17,Part1!$C$1:$C$17,1)'
